$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I (I0) and J (IF), styled like the other
# header cells (same style as H1 / "IP").
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Per-row values for I (I0) and J (IF) columns.
$values = @{
    2  = @(1, 7)
    3  = @(1, 5)
    4  = @(1, 5)
    5  = @(1, 3)
    6  = @(1, 6)
    7  = @(1, 6)
    8  = @(1, 7)
    9  = @(1, 8)
    10 = @(1, 6)
    11 = @(1, 6)
    12 = @(1, 5)
    13 = @(1, 6)
    14 = @(1, 6)
    15 = @(1, 5)
    16 = @(1, 5)
    17 = @(1, 4)
    18 = @(1, 5)
    19 = @(1, 5)
    20 = @(3, 6)
    21 = @(1, 3)
    22 = @(1, 2)
}

foreach ($r in $values.Keys) {
    $pair = $values[$r]
    $ws.Cells.Item($r, 9).Value = $pair[0]
    $ws.Cells.Item($r, 10).Value = $pair[1]
}
